$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 474.52
$ws.Range("I33").Value = 109.166664
$ws.Range("J33").Value = 1414
$ws.Range("K33").Value = 109.166664
$ws.Range("L33").Value = 1414
$ws.Range("M33").Value = 119.833336
$ws.Range("N33").Value = -1872

$ws.Range("H116").Value = 2040
$ws.Range("J116").Value = 2600
$ws.Range("L116").Value = 2600
$ws.Range("N116").Value = -9484

$ws.Range("H132").Value = 1402764.6
$ws.Range("I132").Value = 2612.5173
$ws.Range("J132").Value = 8170166.5
$ws.Range("K132").Value = 7837.5519
$ws.Range("L132").Value = 24510499.5
$ws.Range("M132").Value = -5307.5519
$ws.Range("N132").Value = -24515559.5

$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 20000
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21466

$ws.Range("H74").Value = 10501893
$ws.Range("I74").Value = 13212065
$ws.Range("J74").Value = 203240
$ws.Range("K74").Value = 13212065
$ws.Range("L74").Value = 203240
$ws.Range("M74").Value = -13211191
$ws.Range("N74").Value = -204988

$ws.Range("H77").Value = 10501893
$ws.Range("I77").Value = 13212065
$ws.Range("J77").Value = 203240
$ws.Range("K77").Value = 66060325
$ws.Range("L77").Value = 1016200
$ws.Range("M77").Value = -66055957
$ws.Range("N77").Value = -1024936

$ws.Range("H82").Value = 40181
$ws.Range("J82").Value = 40181
$ws.Range("L82").Value = 40181
$ws.Range("N82").Value = -40903

$ws.Range("H85").Value = 40181
$ws.Range("J85").Value = 40181
$ws.Range("L85").Value = 40181
$ws.Range("N85").Value = -42677

$ws.Range("H122").Value = 5850023
$ws.Range("I122").Value = 2566.8
$ws.Range("J122").Value = 12347197
$ws.Range("K122").Value = 7700.400000000001
$ws.Range("L122").Value = 37041591
$ws.Range("M122").Value = -5250.400000000001
$ws.Range("N122").Value = -37046491

$ws.Range("H128").Value = 47999.668
$ws.Range("J128").Value = 47999.668
$ws.Range("L128").Value = 47999.668
$ws.Range("N128").Value = -57959.668

$ws.Range("H134").Value = 50314.5
$ws.Range("J134").Value = 50314.5
$ws.Range("L134").Value = 50314.5
$ws.Range("N134").Value = -60454.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3622.9355
$ws.Range("I31").Value = 1312.44
$ws.Range("J31").Value = 13250
$ws.Range("K31").Value = 1312.44
$ws.Range("L31").Value = 13250
$ws.Range("M31").Value = -1017.44
$ws.Range("N31").Value = -13840

$ws.Range("H32").Value = 2280
$ws.Range("I32").Value = 2100
$ws.Range("K32").Value = 2100
$ws.Range("M32").Value = -1784

$ws.Range("H34").Value = 3622.9355
$ws.Range("I34").Value = 1312.44
$ws.Range("J34").Value = 13250
$ws.Range("K34").Value = 1312.44
$ws.Range("L34").Value = 13250
$ws.Range("M34").Value = -1110.44
$ws.Range("N34").Value = -13654

$ws.Range("H64").Value = 21110
$ws.Range("J64").Value = 21110
$ws.Range("L64").Value = 21110
$ws.Range("N64").Value = -21606

$ws.Range("H67").Value = 21110
$ws.Range("J67").Value = 21110
$ws.Range("L67").Value = 21110
$ws.Range("N67").Value = -22826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1084.826
$ws.Range("I131").Value = 512.9
$ws.Range("J131").Value = 1243.6945
$ws.Range("K131").Value = 1538.7
$ws.Range("L131").Value = 3731.0835
$ws.Range("M131").Value = 3501.3
$ws.Range("N131").Value = -13811.0835

$ws.Range("H137").Value = 44142.715
$ws.Range("I137").Value = 999.6
$ws.Range("J137").Value = 68111.11
$ws.Range("K137").Value = 2998.8
$ws.Range("L137").Value = 204333.33
$ws.Range("M137").Value = 2101.2
$ws.Range("N137").Value = -214533.33

$ws.Range("H138").Value = 3406.842
$ws.Range("I138").Value = 2647.7778
$ws.Range("J138").Value = 4090
$ws.Range("K138").Value = 7943.3334
$ws.Range("L138").Value = 12270
$ws.Range("M138").Value = -2803.3334
$ws.Range("N138").Value = -22550

$ws.Range("H139").Value = 4215.154
$ws.Range("I139").Value = 2267.7896
$ws.Range("J139").Value = 5336.364
$ws.Range("K139").Value = 6803.3688
$ws.Range("L139").Value = 16009.092
$ws.Range("M139").Value = -1663.3688
$ws.Range("N139").Value = -26289.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 766.6667
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 650
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 650
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -3894

$ws.Range("H122").Value = 1166.6666
$ws.Range("I122").Value = 1166.6666
$ws.Range("K122").Value = 3499.9998
$ws.Range("M122").Value = -1049.9998

$ws.Range("H132").Value = 58671.57
$ws.Range("I132").Value = 34790.168
$ws.Range("J132").Value = 201960
$ws.Range("K132").Value = 104370.504
$ws.Range("L132").Value = 605880
$ws.Range("M132").Value = -101840.504
$ws.Range("N132").Value = -610940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 31345.5
$ws.Range("J64").Value = 31345.5
$ws.Range("L64").Value = 31345.5
$ws.Range("N64").Value = -31795.5

$ws.Range("H67").Value = 31345.5
$ws.Range("J67").Value = 31345.5
$ws.Range("L67").Value = 31345.5
$ws.Range("N67").Value = -32905.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 36107.25
$ws.Range("J46").Value = 42143
$ws.Range("L46").Value = 42143
$ws.Range("N46").Value = -42605

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = $null

$ws.Range("H134").Value = 36107.25
$ws.Range("J134").Value = 42143
$ws.Range("L134").Value = 126429
$ws.Range("N134").Value = -131499
